# Applies the betexplorer.com Premier League 2023-2024 update:
#  - re-orders the match rows in a few places (rows 78/80, 83-87,
#    94/95, 98/100) -- the "index"(A) and "data_partida"(E) columns
#    stay put, only the match detail columns F:V move between rows;
#  - appends one new match row (Chelsea vs Manchester City, index 120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow($row, $vals) {
    # $vals holds F..V (17 values) for $row, in column order.
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($row, 6 + $i).Value = $vals[$i]
    }
}

Set-MatchRow 78 @('Brighton', 2, 'Liverpool', 2, 2.47, '24/09/2023 10:01', 3.26, '08/10/2023 14:44', 3.9, '24/09/2023 10:01', 4.19, '08/10/2023 14:44', 2.54, '24/09/2023 10:01', 2.1, '08/10/2023 14:44', 'https://www.betexplorer.com/football/england/premier-league/brighton-liverpool/2m5wFPdk/')
Set-MatchRow 80 @('Wolves', 1, 'Aston Villa', 1, 2.81, '24/09/2023 10:02', 3.7, '08/10/2023 14:45', 3.39, '24/09/2023 10:02', 3.68, '08/10/2023 14:45', 2.62, '24/09/2023 10:02', 2.08, '08/10/2023 14:45', 'https://www.betexplorer.com/football/england/premier-league/wolves-aston-villa/GAT6GxYg/')
Set-MatchRow 83 @('Nottingham', 2, 'Luton', 2, 1.78, '02/10/2023 08:30', 1.76, '21/10/2023 15:56', 3.63, '02/10/2023 08:30', 3.65, '21/10/2023 15:58', 4.53, '02/10/2023 08:30', 5.39, '21/10/2023 15:58', 'https://www.betexplorer.com/football/england/premier-league/nottingham-luton/tC3uVymm/')
Set-MatchRow 84 @('Newcastle', 4, 'Crystal Palace', 0, 1.49, '01/10/2023 23:01', 1.48, '21/10/2023 15:50', 4.51, '01/10/2023 23:01', 4.49, '21/10/2023 15:58', 7.02, '01/10/2023 23:01', 7.68, '21/10/2023 15:58', 'https://www.betexplorer.com/football/england/premier-league/newcastle-utd-crystal-palace/2L4yWHXt/')
Set-MatchRow 85 @('Manchester City', 2, 'Brighton', 1, 1.47, '01/10/2023 23:01', 1.36, '21/10/2023 15:34', 4.98, '01/10/2023 23:01', 5.69, '21/10/2023 15:58', 6.36, '01/10/2023 23:01', 7.74, '21/10/2023 15:58', 'https://www.betexplorer.com/football/england/premier-league/manchester-city-brighton/ptI9zbPP/')
Set-MatchRow 86 @('Bournemouth', 1, 'Wolves', 2, 2.24, '01/10/2023 23:01', 2.31, '21/10/2023 15:55', 3.47, '01/10/2023 23:01', 3.46, '21/10/2023 15:57', 3.33, '01/10/2023 23:01', 3.32, '21/10/2023 15:58', 'https://www.betexplorer.com/football/england/premier-league/bournemouth-wolves/bZIBFdm0/')
Set-MatchRow 87 @('Brentford', 3, 'Burnley', 0, 1.63, '02/10/2023 08:30', 1.74, '21/10/2023 15:45', 3.97, '02/10/2023 08:30', 3.9, '21/10/2023 15:58', 5.07, '02/10/2023 08:30', 5.09, '21/10/2023 15:58', 'https://www.betexplorer.com/football/england/premier-league/brentford-burnley/6aMJDzIC/')
Set-MatchRow 94 @('Bournemouth', 2, 'Burnley', 1, 2.22, '10/10/2023 14:20', 2.13, '28/10/2023 15:57', 3.45, '10/10/2023 14:20', 3.65, '28/10/2023 15:59', 3.18, '10/10/2023 14:20', 3.51, '28/10/2023 15:57', 'https://www.betexplorer.com/football/england/premier-league/bournemouth-burnley/W0dhSZW5/')
Set-MatchRow 95 @('Arsenal', 5, 'Sheffield Utd', 0, 1.17, '10/10/2023 14:22', 1.13, '28/10/2023 14:59', 7.34, '10/10/2023 14:22', 9.5, '28/10/2023 14:59', 12.88, '10/10/2023 14:22', 21, '28/10/2023 14:59', 'https://www.betexplorer.com/football/england/premier-league/arsenal-sheffield-utd/p40dRgnC/')
Set-MatchRow 98 @('Liverpool', 3, 'Nottingham', 0, 1.23, '10/10/2023 14:02', 1.22, '29/10/2023 14:51', 6.65, '10/10/2023 14:02', 7.31, '29/10/2023 14:59', 9.43, '10/10/2023 14:02', 12.08, '29/10/2023 14:59', 'https://www.betexplorer.com/football/england/premier-league/liverpool-nottingham/IcEJreHn/')
Set-MatchRow 100 @('Aston Villa', 3, 'Luton', 1, 1.38, '10/10/2023 14:32', 1.34, '29/10/2023 14:55', 4.92, '10/10/2023 14:32', 5.85, '29/10/2023 14:55', 7.42, '10/10/2023 14:32', 8.55, '29/10/2023 14:57', 'https://www.betexplorer.com/football/england/premier-league/aston-villa-luton/SSk1QD1I/')

# New last row (index 120 / sheet row 121) -- clone row 120's cell
# formatting (bold+bordered index cell, datetime-formatted date cell)
# onto row 121, then overwrite every value for the new match.
$ws.Range("A120:V120").Copy($ws.Range("A121:V121"))

$ws.Range("A121").Value = 120
$ws.Range("B121").Value = 'england'
$ws.Range("C121").Value = 'premier-league'
$ws.Range("D121").Value = '2023-2024'
$ws.Range("E121").Value = 45242.72916666666
$ws.Range("F121").Value = 'Chelsea'
$ws.Range("G121").Value = 4
$ws.Range("H121").Value = 'Manchester City'
$ws.Range("I121").Value = 4
$ws.Range("J121").Value = 3.78
$ws.Range("K121").Value = '28/10/2023 22:02'
$ws.Range("L121").Value = 4.8
$ws.Range("M121").Value = '12/11/2023 17:17'
$ws.Range("N121").Value = 3.64
$ws.Range("O121").Value = '28/10/2023 22:02'
$ws.Range("P121").Value = 4
$ws.Range("Q121").Value = '12/11/2023 17:17'
$ws.Range("R121").Value = 2
$ws.Range("S121").Value = '28/10/2023 22:02'
$ws.Range("T121").Value = 1.76
$ws.Range("U121").Value = '12/11/2023 17:17'
$ws.Range("V121").Value = 'https://www.betexplorer.com/football/england/premier-league/chelsea-manchester-city/UVvS9XFo/'
